$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name/Link swap for rows 10 and 11 (BinanceUSD <-> OKB) ---
$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(11, 2).Value = "BinanceUSD"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"

# --- Price (column D) updates; force text so numeric-looking strings are not
#     auto-converted to numbers (matches the original inlineStr text cells) ---
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "24.790.50"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.705.40"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "315.09"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4003"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4047"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "1.472"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "53.63"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.9987"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08805"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "26.20"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.513"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.991"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.00001342"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.720.78"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "95.47"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.07163"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "20.88"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.290"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.002"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "14.46"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "24.793.75"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.885"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "23.12"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.161"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "161.25"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "144.17"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.209"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.278"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.906.86"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08656"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.03191"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "7.282"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.028"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.2847"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.8393"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.09473"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "10.71"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "14.21"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.480"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "17.45"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.717"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.7429"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.216"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.377"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.001"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "140.35"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.08397"

# --- Volume(1h) (column E) updates ---
$ws.Cells.Item(2, 5).Value = "  +0.74%  "
$ws.Cells.Item(3, 5).Value = "  +0.94%  "
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 5).Value = "  +0.32%  "
$ws.Cells.Item(6, 5).Value = "  -0.25%  "
$ws.Cells.Item(7, 5).Value = "  +2.79%  "
$ws.Cells.Item(8, 5).Value = "  +0.31%  "
$ws.Cells.Item(9, 5).Value = "  -1.65%  "
$ws.Cells.Item(10, 5).Value = "  +1.44%  "
$ws.Cells.Item(11, 5).Value = "  -0.48%  "
$ws.Cells.Item(12, 5).Value = "  +0.53%  "
$ws.Cells.Item(13, 5).Value = "  +5.39%  "
$ws.Cells.Item(14, 5).Value = "  -0.75%  "
$ws.Cells.Item(15, 5).Value = "  +0.38%  "
$ws.Cells.Item(16, 5).Value = "  -0.53%  "
$ws.Cells.Item(17, 5).Value = "  +2.33%  "
$ws.Cells.Item(18, 5).Value = "  -3.18%  "
$ws.Cells.Item(19, 5).Value = "  +0.78%  "
$ws.Cells.Item(20, 5).Value = "  +5.09%  "
$ws.Cells.Item(21, 5).Value = "  +0.02%  "
$ws.Cells.Item(22, 5).Value = "  -0.12%  "
$ws.Cells.Item(23, 5).Value = "  +1.26%  "
$ws.Cells.Item(24, 5).Value = "  +0.77%  "
$ws.Cells.Item(25, 5).Value = "  +0.18%  "
$ws.Cells.Item(26, 5).Value = "  -4.21%  "
$ws.Cells.Item(27, 5).Value = "  +1.31%  "
$ws.Cells.Item(28, 5).Value = "  +18.07%  "
$ws.Cells.Item(29, 5).Value = "  -0.50%  "
$ws.Cells.Item(30, 5).Value = "  +5.26%  "
$ws.Cells.Item(31, 5).Value = "  -7.05%  "
$ws.Cells.Item(32, 5).Value = "  +14.97%  "
$ws.Cells.Item(33, 5).Value = "  +2.16%  "
$ws.Cells.Item(34, 5).Value = "  -2.33%  "
$ws.Cells.Item(35, 5).Value = "  +9.31%  "
$ws.Cells.Item(36, 5).Value = "  -2.13%  "
$ws.Cells.Item(37, 5).Value = "  -1.16%  "
$ws.Cells.Item(38, 5).Value = "  +4.01%  "
$ws.Cells.Item(39, 5).Value = "  +7.30%  "
$ws.Cells.Item(40, 5).Value = "  +3.64%  "
$ws.Cells.Item(41, 5).Value = "  -0.62%  "
$ws.Cells.Item(42, 5).Value = "  +0.08%  "
$ws.Cells.Item(43, 5).Value = "  +1.39%  "
$ws.Cells.Item(44, 5).Value = "  +4.65%  "
$ws.Cells.Item(45, 5).Value = "  +4.81%  "
$ws.Cells.Item(46, 5).Value = "  +3.28%  "
$ws.Cells.Item(47, 5).Value = "  +0.54%  "
$ws.Cells.Item(48, 5).Value = "  +3.03%  "
$ws.Cells.Item(49, 5).Value = "  -0.20%  "
$ws.Cells.Item(50, 5).Value = "  +1.73%  "
$ws.Cells.Item(51, 5).Value = "  +5.36%  "
